# Update the "State" column (F) for every task row on the sheet so that
# every entry reflects that the work item is complete ("Finished")
# instead of the previous "Processing" / "UnFinished" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2:F10").Value = "Finished"
